# Risk Score related changes
# - API Type (column D) for rows 14-21 switches from "Internal" to "External"
# - Veracode SLA Breach (L) / Pen Test SLA Breach (M) are (re)populated for every
#   data row: even rows -> "Withen SLA", odd rows -> "SLA Breached"
# - Selection moves to Q12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 14-21: API Type becomes "External"
for ($r = 14; $r -le 21; $r++) {
    $ws.Cells.Item($r, 4).Value = "External"
}

# Rows 2-21: Veracode SLA Breach (L) / Pen Test SLA Breach (M)
for ($r = 2; $r -le 21; $r++) {
    if ($r % 2 -eq 0) {
        $slaValue = "Withen SLA"
    } else {
        $slaValue = "SLA Breached"
    }
    $ws.Cells.Item($r, 12).Value = $slaValue
    $ws.Cells.Item($r, 13).Value = $slaValue
}

# Update the active selection to Q12
$ws.Range("Q12").Select() | Out-Null
